# "Getting Cell Info from Sheet" - append the new "Sushi" purchase row
# (row 8) to Sheet1: a date/time in column A, the item name in column B
# (backed by a new shared string), and the amount in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: purchase date/time, formatted like the other rows above it
# (built-in date-time number format "m/d/yy h:mm" == numFmtId 22).
$ws.Range("A8").Value = 43848.540972222225
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"

# Column B: item name - matches the vertical-center alignment used by
# the other item cells in this column.
$ws.Range("B8").Value = "Sushi"
$ws.Range("B8").VerticalAlignment = -4108

# Column C: amount - matches the right-aligned / vertical-center
# alignment used by the other amount cells in this column.
$ws.Range("C8").Value = 129
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("C8").HorizontalAlignment = -4152

# Selection ends up on C14 after entering the new row of data.
$ws.Range("C14").Select() | Out-Null
